# Add a new verb row ("recruits" / "recruit" / MatchCase="N") into the
# alphabetically-sorted verb table on Sheet1. It belongs between "provides"
# (row 58) and "replicates" (old row 59), so insert a new row 59 and shift
# everything below down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 59 (pushes "replicates" and everything after down by one).
$ws.Rows.Item(59).Insert()

# Column A/B/C literal values. The sheet's convention stores the literal
# quote characters as part of the string (e.g. A2 = the 12-character string
# "accomplishes" including the quote marks), since the columns feed VBA
# source-text generator formulas in D:G.
$ws.Range("A59").Value = '"recruits"'
$ws.Range("B59").Value = '"recruit"'
$ws.Range("C59").Value = '"N"'

# D:G recreate the same generator formulas used by every other data row,
# anchored to row 59 so they participate in the existing shared-formula
# groups exactly like a native Excel row-insert would.
$ws.Range("D59").Formula = '="verbTense("&ROW(A59)-1&", "&1&") = "&A59'
$ws.Range("E59").Formula = '="verbTense("&ROW(A59)-1&", "&2&") = "&B59'
$ws.Range("F59").Formula = '="verbTense("&ROW(A59)-1&", "&3&") = "&C59'
$ws.Range("G59").Formula = '=D59&" : "&E59&" : "&F59'

# Re-apply the table's sort over the (now one-row-larger) data range so the
# sheet's remembered auto-filter/sort state covers the new last row.
$sortRange = $ws.Range("A2:C73")
$ws.Sort.SortFields.Clear()
$ws.Sort.SortFields.Add($ws.Range("A2"))
$ws.Sort.SetRange($sortRange)
$ws.Sort.Header = 2
$ws.Sort.Apply()

# Leave the cursor where the author's saved view shows it.
$ws.Range("B60").Select()
